$wb = $excel.ActiveWorkbook

function Set-Cell {
    param($ws, [string]$addr, [double]$val)
    $ws.Range($addr).Value = $val
}

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
Set-Cell $ws "H21" 50000
Set-Cell $ws "J21" 50000
Set-Cell $ws "L21" 50000
Set-Cell $ws "N21" -50936
Set-Cell $ws "H23" 50000
Set-Cell $ws "J23" 50000
Set-Cell $ws "L23" 50000
Set-Cell $ws "N23" -50468
Set-Cell $ws "H54" 4846
Set-Cell $ws "I54" 2292
Set-Cell $ws "J54" 7400
Set-Cell $ws "K54" 2292
Set-Cell $ws "L54" 7400
Set-Cell $ws "M54" -1806
Set-Cell $ws "N54" -8372
Set-Cell $ws "H125" 10192439
Set-Cell $ws "I125" 619.8
Set-Cell $ws "J125" 18685622
Set-Cell $ws "K125" 5578.2
Set-Cell $ws "L125" 168170598
Set-Cell $ws "M125" -3118.2
Set-Cell $ws "N125" -168175518
Set-Cell $ws "H132" 383127.22
Set-Cell $ws "I132" 405168.84
Set-Cell $ws "K132" 1215506.52
Set-Cell $ws "M132" -1212976.52
Set-Cell $ws "H135" 2487.7144
Set-Cell $ws "I135" 2525.2307
Set-Cell $ws "J135" 2000
Set-Cell $ws "K135" 22727.0763
Set-Cell $ws "L135" 18000
Set-Cell $ws "M135" -20192.0763
Set-Cell $ws "N135" -23070
Set-Cell $ws "H137" 125001920
Set-Cell $ws "I137" 333334340
Set-Cell $ws "J137" 2480
Set-Cell $ws "K137" 1000003020
Set-Cell $ws "L137" 7440
Set-Cell $ws "M137" -1000000470
Set-Cell $ws "N137" -12540

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
Set-Cell $ws "H32" 25480.729
Set-Cell $ws "I32" 4142.1665
Set-Cell $ws "K32" 4142.1665
Set-Cell $ws "M32" -3855.1665
Set-Cell $ws "H45" 1191.6923
Set-Cell $ws "I45" 1099.2
Set-Cell $ws "K45" 1099.2
Set-Cell $ws "M45" -722.2
Set-Cell $ws "H63" 21286.428
Set-Cell $ws "I63" 26801
Set-Cell $ws "K63" 26801
Set-Cell $ws "M63" -26115
Set-Cell $ws "H66" 21286.428
Set-Cell $ws "I66" 26801
Set-Cell $ws "K66" 134005
Set-Cell $ws "M66" -130573
Set-Cell $ws "H122" 9681
Set-Cell $ws "I122" 2908
Set-Cell $ws "J122" 30000
Set-Cell $ws "K122" 8724
Set-Cell $ws "L122" 90000
Set-Cell $ws "M122" -6274
Set-Cell $ws "N122" -94900
Set-Cell $ws "H123" 41689.855
Set-Cell $ws "J123" 41689.855
Set-Cell $ws "L123" 41689.855
Set-Cell $ws "N123" -51489.855
Set-Cell $ws "H141" 36000
Set-Cell $ws "J141" 36000
Set-Cell $ws "L141" 36000
Set-Cell $ws "N141" -46360

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
Set-Cell $ws "H105" 266112.28
Set-Cell $ws "I105" 2793.077
Set-Cell $ws "J105" 836637.25
Set-Cell $ws "K105" 2793.077
Set-Cell $ws "L105" 836637.25
Set-Cell $ws "M105" -1046.077
Set-Cell $ws "N105" -840131.25
Set-Cell $ws "H107" 1375.125
Set-Cell $ws "I107" 1375.125
Set-Cell $ws "K107" 1375.125
Set-Cell $ws "M107" 544.875

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
Set-Cell $ws "H31" 3179.9
Set-Cell $ws "I31" 1140
Set-Cell $ws "J31" 5847.4614
Set-Cell $ws "K31" 1140
Set-Cell $ws "L31" 5847.4614
Set-Cell $ws "M31" -845
Set-Cell $ws "N31" -6437.4614
Set-Cell $ws "H34" 3179.9
Set-Cell $ws "I34" 1140
Set-Cell $ws "J34" 5847.4614
Set-Cell $ws "K34" 1140
Set-Cell $ws "L34" 5847.4614
Set-Cell $ws "M34" -938
Set-Cell $ws "N34" -6251.4614
Set-Cell $ws "H99" 5216740.5
Set-Cell $ws "J99" 2000
Set-Cell $ws "L99" 2000
Set-Cell $ws "N99" -4996
Set-Cell $ws "H126" 5216740.5
Set-Cell $ws "J126" 2000
Set-Cell $ws "L126" 6000
Set-Cell $ws "N126" -10940
Set-Cell $ws "H132" 2920.3635
Set-Cell $ws "I132" 2412.7693
Set-Cell $ws "J132" 4805.7144
Set-Cell $ws "K132" 7238.3079
Set-Cell $ws "L132" 14417.1432
Set-Cell $ws "M132" -4708.3079
Set-Cell $ws "N132" -19477.1432
Set-Cell $ws "H134" 2855.361
Set-Cell $ws "I134" 1271
Set-Cell $ws "J134" 5345.0713
Set-Cell $ws "K134" 3813
Set-Cell $ws "L134" 16035.2139
Set-Cell $ws "M134" -1278
Set-Cell $ws "N134" -21105.2139

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
Set-Cell $ws "H4" 12603.708
Set-Cell $ws "I4" 96.611115
Set-Cell $ws "J4" 50125
Set-Cell $ws "K4" 289.833345
Set-Cell $ws "L4" 150375
Set-Cell $ws "M4" -177.833345
Set-Cell $ws "N4" -150599
Set-Cell $ws "H16" 962.625
Set-Cell $ws "I16" 850.5
Set-Cell $ws "J16" 1000
Set-Cell $ws "K16" 2551.5
Set-Cell $ws "L16" 3000
Set-Cell $ws "M16" -2378.5
Set-Cell $ws "N16" -3346
Set-Cell $ws "H54" 4900
Set-Cell $ws "J54" 4900
Set-Cell $ws "L54" 14700
Set-Cell $ws "N54" -15818
Set-Cell $ws "H134" 5915.143
Set-Cell $ws "I134" 3645.7778
Set-Cell $ws "J134" 10000
Set-Cell $ws "K134" 10937.3334
Set-Cell $ws "L134" 30000
Set-Cell $ws "M134" -5867.3334
Set-Cell $ws "N134" -40140

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
Set-Cell $ws "H102" 3582.9092
Set-Cell $ws "I102" 3626.5
Set-Cell $ws "J102" 3466.6667
Set-Cell $ws "K102" 3626.5
Set-Cell $ws "L102" 3466.6667
Set-Cell $ws "M102" -2004.5
Set-Cell $ws "N102" -6710.6667
Set-Cell $ws "H122" 742835.9399999999
Set-Cell $ws "I122" 1235877
Set-Cell $ws "J122" 3274.3333
Set-Cell $ws "K122" 3707631
Set-Cell $ws "L122" 9822.999899999999
Set-Cell $ws "M122" -3705181
Set-Cell $ws "N122" -14722.9999
Set-Cell $ws "H123" 12138.167
Set-Cell $ws "J123" 12138.167
Set-Cell $ws "L123" 12138.167
Set-Cell $ws "N123" -17038.167
Set-Cell $ws "H132" 3382.0527
Set-Cell $ws "I132" 3230.8
Set-Cell $ws "K132" 9692.400000000001
Set-Cell $ws "M132" -7162.400000000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
Set-Cell $ws "H40" 3005.4736
Set-Cell $ws "I40" 2001.3334
Set-Cell $ws "J40" 3193.75
Set-Cell $ws "K40" 2001.3334
Set-Cell $ws "L40" 3193.75
Set-Cell $ws "M40" -1865.3334
Set-Cell $ws "N40" -3465.75
Set-Cell $ws "H122" 3773.8262
Set-Cell $ws "I122" 3318
Set-Cell $ws "J122" 3934.7058
Set-Cell $ws "K122" 9954
Set-Cell $ws "L122" 11804.1174
Set-Cell $ws "M122" -7504
Set-Cell $ws "N122" -16704.1174

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
Set-Cell $ws "H122" 101278.8
Set-Cell $ws "I122" 126036
Set-Cell $ws "K122" 378108
Set-Cell $ws "M122" -375658
Set-Cell $ws "H123" 30727.273
Set-Cell $ws "J123" 30727.273
Set-Cell $ws "L123" 30727.273
Set-Cell $ws "N123" -40527.273
Set-Cell $ws "H132" 3251.3157
Set-Cell $ws "I132" 3019.7856
Set-Cell $ws "J132" 3899.6
Set-Cell $ws "K132" 9059.356800000001
Set-Cell $ws "L132" 11698.8
Set-Cell $ws "M132" -6529.356800000001
Set-Cell $ws "N132" -16758.8
